# Apply the changes described by the commit "updated 4.0 files and mdl"
# to the Fuel Prod Imp Exp Balancing Priorities workbook.

$wb = $excel.ActiveWorkbook

# --- About sheet -----------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

# "Last updated" date moves from 1/3/2024 to 3/28/2024 (serial 45294 -> 45379)
$wsAbout.Range("C1").Value = 45379

# --- FPIEBP sheet ------------------------------------------------------
$wsFPIEBP = $wb.Worksheets.Item("FPIEBP")

# Re-prioritize "hard coal" row: production/imports/exports priority values
# B3 (production): 3 -> 1
# C3 (imports):    2 -> 3
# D3 (exports):    1 -> 2
$wsFPIEBP.Range("B3").Value = 1
$wsFPIEBP.Range("C3").Value = 3
$wsFPIEBP.Range("D3").Value = 2

# Update the active selection on the FPIEBP sheet to E3 (from F4)
$wsFPIEBP.Activate() | Out-Null
$wsFPIEBP.Range("E3").Select() | Out-Null
